$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.382.52"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "1.878.84"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7163"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07976"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3144"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08080"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.71%  "
$ws.Range("D12").Value = "1.874.83"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.79%  "
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7080"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.390"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008432"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "29.384.30"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").Value = "2.133.68"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.676"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1578"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.065"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.418"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.315"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.227"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.943"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7569"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.175"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.704"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "1.282.27"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.765"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.405"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("E42").Value = "  +1.56%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "111.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "74.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("E46").Value = "  -1.11%  "
$ws.Range("D47").Value = "2.027.63"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5207"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.531"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4341"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.44%  "
